# Auto-generated Excel COM-interop script to apply scheduled runner price updates
# across all 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1349.76
$ws.Range("I28").Value = 1119.8334
$ws.Range("K28").Value = 1119.8334
$ws.Range("M28").Value = -634.8334
$ws.Range("H53").Value = 258
$ws.Range("I53").Value = 144.5
$ws.Range("J53").Value = 333.66666
$ws.Range("K53").Value = 144.5
$ws.Range("L53").Value = 333.66666
$ws.Range("M53").Value = 492.5
$ws.Range("N53").Value = -1607.66666
$ws.Range("H76").Value = 7206.375
$ws.Range("I76").Value = 9995.666999999999
$ws.Range("J76").Value = 5532.8
$ws.Range("K76").Value = 9995.666999999999
$ws.Range("L76").Value = 5532.8
$ws.Range("M76").Value = -9680.666999999999
$ws.Range("N76").Value = -6162.8
$ws.Range("H79").Value = 7206.375
$ws.Range("I79").Value = 9995.666999999999
$ws.Range("J79").Value = 5532.8
$ws.Range("K79").Value = 9995.666999999999
$ws.Range("L79").Value = 5532.8
$ws.Range("M79").Value = -8903.666999999999
$ws.Range("N79").Value = -7716.8
$ws.Range("H98").Value = 653.2778
$ws.Range("I98").Value = 653.2778
$ws.Range("K98").Value = 653.2778
$ws.Range("M98").Value = 844.7222
$ws.Range("H106").Value = 1198.6
$ws.Range("I106").Value = 991.53845
$ws.Range("K106").Value = 991.53845
$ws.Range("M106").Value = -360.53845
$ws.Range("H107").Value = 753.2
$ws.Range("I107").Value = 753.2
$ws.Range("K107").Value = 753.2
$ws.Range("M107").Value = 1166.8
$ws.Range("H118").Value = 77895.08
$ws.Range("I118").Value = 77895.08
$ws.Range("K118").Value = 233685.24
$ws.Range("M118").Value = -232028.24
$ws.Range("H122").Value = 653.2778
$ws.Range("I122").Value = 653.2778
$ws.Range("K122").Value = 1959.8334
$ws.Range("M122").Value = 490.1666

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2243.4707
$ws.Range("I2").Value = 2327.8823
$ws.Range("J2").Value = 2159.0588
$ws.Range("K2").Value = 2327.8823
$ws.Range("L2").Value = 2159.0588
$ws.Range("M2").Value = -2214.8823
$ws.Range("N2").Value = -2385.0588
$ws.Range("H61").Value = 4341.231
$ws.Range("J61").Value = 5813
$ws.Range("L61").Value = 5813
$ws.Range("N61").Value = -6237
$ws.Range("H97").Value = 2883.3333
$ws.Range("I97").Value = 3400.25
$ws.Range("J97").Value = 1849.5
$ws.Range("K97").Value = 3400.25
$ws.Range("L97").Value = 1849.5
$ws.Range("M97").Value = -2904.25
$ws.Range("N97").Value = -2841.5
$ws.Range("H104").Value = 96237.836
$ws.Range("J104").Value = 115285.4
$ws.Range("L104").Value = 115285.4
$ws.Range("N104").Value = -122273.4
$ws.Range("H116").Value = 2243.4707
$ws.Range("I116").Value = 2327.8823
$ws.Range("J116").Value = 2159.0588
$ws.Range("K116").Value = 2327.8823
$ws.Range("L116").Value = 2159.0588
$ws.Range("M116").Value = -33.88230000000021
$ws.Range("N116").Value = -6747.0588
$ws.Range("H122").Value = 1672.0454
$ws.Range("I122").Value = 1304.5264
$ws.Range("K122").Value = 3913.5792
$ws.Range("M122").Value = -1463.5792
$ws.Range("H136").Value = 4341.231
$ws.Range("J136").Value = 5813
$ws.Range("L136").Value = 17439
$ws.Range("N136").Value = -22539

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2243.4707
$ws.Range("I3").Value = 2327.8823
$ws.Range("J3").Value = 2159.0588
$ws.Range("K3").Value = 2327.8823
$ws.Range("L3").Value = 2159.0588
$ws.Range("M3").Value = -2213.8823
$ws.Range("N3").Value = -2387.0588
$ws.Range("H94").Value = 1393.5
$ws.Range("I94").Value = 1176.7742
$ws.Range("J94").Value = 2737.2
$ws.Range("K94").Value = 1176.7742
$ws.Range("L94").Value = 2737.2
$ws.Range("M94").Value = -725.7742000000001
$ws.Range("N94").Value = -3639.2
$ws.Range("H107").Value = 1170.1666
$ws.Range("J107").Value = 1455.875
$ws.Range("L107").Value = 1455.875
$ws.Range("N107").Value = -5295.875
$ws.Range("H132").Value = 86375
$ws.Range("J132").Value = 86375
$ws.Range("L132").Value = 86375
$ws.Range("N132").Value = -96495
$ws.Range("H138").Value = 20000
$ws.Range("I138").Value = 20000
$ws.Range("K138").Value = 20000
$ws.Range("M138").Value = -14860
$ws.Range("H139").Value = 82161
$ws.Range("J139").Value = 82701.25
$ws.Range("L139").Value = 82701.25
$ws.Range("N139").Value = -92981.25

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2249.75
$ws.Range("I16").Value = 2249.75
$ws.Range("K16").Value = 2249.75
$ws.Range("M16").Value = -1962.75
$ws.Range("H99").Value = 2943.2222
$ws.Range("I99").Value = 2498.4285
$ws.Range("J99").Value = 4500
$ws.Range("K99").Value = 2498.4285
$ws.Range("L99").Value = 4500
$ws.Range("M99").Value = -1000.4285
$ws.Range("N99").Value = -7496
$ws.Range("H105").Value = 1827.125
$ws.Range("I105").Value = 1611.1666
$ws.Range("K105").Value = 1611.1666
$ws.Range("M105").Value = 135.8334
$ws.Range("H107").Value = 341.375
$ws.Range("I107").Value = 349.07144
$ws.Range("J107").Value = 287.5
$ws.Range("K107").Value = 349.07144
$ws.Range("L107").Value = 287.5
$ws.Range("M107").Value = 1570.92856
$ws.Range("N107").Value = -4127.5
$ws.Range("H113").Value = 2249.75
$ws.Range("I113").Value = 2249.75
$ws.Range("K113").Value = 2249.75
$ws.Range("M113").Value = -79.75
$ws.Range("H126").Value = 2943.2222
$ws.Range("I126").Value = 2498.4285
$ws.Range("J126").Value = 4500
$ws.Range("K126").Value = 7495.2855
$ws.Range("L126").Value = 13500
$ws.Range("M126").Value = -5025.2855
$ws.Range("N126").Value = -18440
$ws.Range("H135").Value = 74966.664
$ws.Range("J135").Value = 74966.664
$ws.Range("L135").Value = 74966.664
$ws.Range("N135").Value = -85106.664
$ws.Range("H137").Value = 92985.71000000001
$ws.Range("J137").Value = 92985.71000000001
$ws.Range("L137").Value = 92985.71000000001
$ws.Range("N137").Value = -103185.71

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2986
$ws.Range("J68").Value = 3314.6667
$ws.Range("L68").Value = 9944.000100000001
$ws.Range("N68").Value = -11566.0001
$ws.Range("H71").Value = 2986
$ws.Range("J71").Value = 3314.6667
$ws.Range("L71").Value = 29832.0003
$ws.Range("N71").Value = -37944.0003

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 28799.834
$ws.Range("J46").Value = 38600
$ws.Range("L46").Value = 38600
$ws.Range("N46").Value = -38912
$ws.Range("H80").Value = 8372.333000000001
$ws.Range("I80").Value = 2970
$ws.Range("J80").Value = 11073.5
$ws.Range("K80").Value = 2970
$ws.Range("L80").Value = 11073.5
$ws.Range("M80").Value = -1972
$ws.Range("N80").Value = -13069.5
$ws.Range("H83").Value = 8372.333000000001
$ws.Range("I83").Value = 2970
$ws.Range("J83").Value = 11073.5
$ws.Range("K83").Value = 14850
$ws.Range("L83").Value = 55367.5
$ws.Range("M83").Value = -9858
$ws.Range("N83").Value = -65351.5
$ws.Range("H102").Value = 2196.2144
$ws.Range("I102").Value = 2212.9614
$ws.Range("J102").Value = 1978.5
$ws.Range("K102").Value = 2212.9614
$ws.Range("L102").Value = 1978.5
$ws.Range("M102").Value = -590.9614000000001
$ws.Range("N102").Value = -5222.5
$ws.Range("H107").Value = 1322
$ws.Range("H122").Value = 1026.7
$ws.Range("I122").Value = 977.6667
$ws.Range("K122").Value = 2933.0001
$ws.Range("M122").Value = -483.0001000000002
$ws.Range("H126").Value = 6819.32
$ws.Range("I126").Value = 7704.3687
$ws.Range("K126").Value = 23113.1061
$ws.Range("M126").Value = -20643.1061

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9536.416999999999
$ws.Range("I7").Value = 10443.2
$ws.Range("J7").Value = 5002.5
$ws.Range("K7").Value = 10443.2
$ws.Range("L7").Value = 5002.5
$ws.Range("M7").Value = -10331.2
$ws.Range("N7").Value = -5226.5
$ws.Range("H40").Value = 5617.4165
$ws.Range("I40").Value = 4990.4443
$ws.Range("K40").Value = 4990.4443
$ws.Range("M40").Value = -4854.4443
$ws.Range("H126").Value = 9536.416999999999
$ws.Range("I126").Value = 10443.2
$ws.Range("J126").Value = 5002.5
$ws.Range("K126").Value = 31329.6
$ws.Range("L126").Value = 15007.5
$ws.Range("M126").Value = -28859.6
$ws.Range("N126").Value = -19947.5
$ws.Range("H132").Value = 4329.9287
$ws.Range("I132").Value = 3137.6365
$ws.Range("K132").Value = 9412.9095
$ws.Range("M132").Value = -6882.9095

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 466.66666
$ws.Range("I107").Value = 493
$ws.Range("J107").Value = 374.5
$ws.Range("K107").Value = 1479
$ws.Range("L107").Value = 1123.5
$ws.Range("M107").Value = 441
$ws.Range("N107").Value = -4963.5
$ws.Range("H122").Value = 2340.1035
$ws.Range("I122").Value = 2322.6538
$ws.Range("J122").Value = 2491.3333
$ws.Range("K122").Value = 6967.9614
$ws.Range("L122").Value = 7473.999899999999
$ws.Range("M122").Value = -4517.9614
$ws.Range("N122").Value = -12373.9999
$ws.Range("H137").Value = 125984
$ws.Range("J137").Value = 125984
$ws.Range("L137").Value = 125984
$ws.Range("N137").Value = -136184
